$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 values (Q1 2021 revisions) ---
$ws.Range("K75").Value = 8841
$ws.Range("M75").Value = 3477
$ws.Range("O75").Value = 5364
$ws.Range("P75").Value = 615
$ws.Range("Q75").Value = 576
$ws.Range("S75").Value = 1422
$ws.Range("V75").Value = 769
$ws.Range("W75").Value = 165235
$ws.Range("X75").Value = 96871
$ws.Range("Y75").Value = 1274
$ws.Range("AB75").Value = 61135
$ws.Range("AC75").Value = 98253
$ws.Range("AG75").Value = 5952
$ws.Range("AH75").Value = 13461
$ws.Range("AJ75").Value = 1643
$ws.Range("AK75").Value = 11819

# --- Append new row 76 (01-04-2021) ---
# Column A holds a date-like label that must stay literal text ("01-04-2021"),
# not get auto-converted to a date serial number. Format the cell as Text
# before assigning the value, then drop back to the workbook's default
# (unstyled) cell style so the new row matches the rest of the sheet.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = 857
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 0
$ws.Range("E76").Value = 857
$ws.Range("F76").Value = 37325
$ws.Range("G76").Value = 29055
$ws.Range("H76").Value = 107
$ws.Range("I76").Value = 1293
$ws.Range("J76").Value = 6870
$ws.Range("K76").Value = 8753
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 3761
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 4992
$ws.Range("P76").Value = 740
$ws.Range("Q76").Value = 701
$ws.Range("R76").Value = 40
$ws.Range("S76").Value = 1405
$ws.Range("T76").Value = 325
$ws.Range("U76").Value = 326
$ws.Range("V76").Value = 754
$ws.Range("W76").Value = 165558
$ws.Range("X76").Value = 97699
$ws.Range("Y76").Value = 1297
$ws.Range("Z76").Value = 5068
$ws.Range("AA76").Value = 913
$ws.Range("AB76").Value = 60581
$ws.Range("AC76").Value = 100614
$ws.Range("AD76").Value = 84017
$ws.Range("AE76").Value = 4671
$ws.Range("AF76").Value = 5216
$ws.Range("AG76").Value = 6711
$ws.Range("AH76").Value = 13643
$ws.Range("AI76").Value = 0
$ws.Range("AJ76").Value = 1716
$ws.Range("AK76").Value = 11928
